$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append one new row (row 3) with the same "everything stored as text"
# shape as the existing rows (A1:H2), including the quantity column
# ("2"), which looks numeric but must stay text like C2's "23".

# Force text storage for the quantity cell so "2" isn't coerced to a
# number.
$ws.Cells.Item(3, 3).NumberFormat = "@"

# A3 mirrors A2, which is an explicit empty *text* value (not a blank
# cell) - a leading apostrophe forces Excel to store an empty string as
# text instead of leaving the cell blank.
$ws.Cells.Item(3, 1).Value = "'"
$ws.Cells.Item(3, 2).Value = "أحمد شريم"
$ws.Cells.Item(3, 3).Value = "2"
$ws.Cells.Item(3, 4).Value = "الصمود"
$ws.Cells.Item(3, 5).Value = "الرحلة 2"
$ws.Cells.Item(3, 6).Value = "C1"
$ws.Cells.Item(3, 7).Value = "UNICEF"
$ws.Cells.Item(3, 8).Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٦:٣٠:٠٠ م"

# The text-forcing tricks above (quote-prefix / custom number format)
# stamp a new cell style; reset each new cell back to row 2's style so
# row 3 doesn't pick up a spurious style index versus the source rows.
for ($c = 1; $c -le 8; $c++) {
  $ws.Cells.Item(3, $c).Style = $ws.Cells.Item(2, $c).Style
}
